# Generate Report for Handback
#
# - Updates the "Status" text (shared across Overview/zh-cn/de-de) from
#   "Ready for handoff" to "Handed back: in sync with en-US".
# - Adds "Latest Target File" / "Latest Handback File" hyperlink cells
#   (columns E/F) on the zh-cn and de-de sheets, mirroring the existing
#   "Source File Name" (A) / "Latest Handoff File" (C) hyperlinks since
#   the handback is in sync with the handoff.
# - Stamps the "Latest Handback DateTime" column (G) with the real
#   handback timestamp instead of the epoch placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$md1 = "39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.md"
$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/d19d66efe81be9f8b7961c2acc10d0bfded39f9b/e2e/39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.md"

$md2 = "54d53a45-672b-45eb-a5f5-567c408ad3a1.md"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/d19d66efe81be9f8b7961c2acc10d0bfded39f9b/e2e/54d53a45-672b-45eb-a5f5-567c408ad3a1.md"

$zhXlf1 = "39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.b0c490ba30254cd1a1156d1031ac43e4c04209be.zh-cn.xlf"
$zhXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24ca814670d5cd44b52ef305e6d881b55fb74050/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.b0c490ba30254cd1a1156d1031ac43e4c04209be.zh-cn.xlf"

$zhXlf2 = "54d53a45-672b-45eb-a5f5-567c408ad3a1.10def667fa8a4b9195e9896236ea5151b05faccf.zh-cn.xlf"
$zhXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/24ca814670d5cd44b52ef305e6d881b55fb74050/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/54d53a45-672b-45eb-a5f5-567c408ad3a1.10def667fa8a4b9195e9896236ea5151b05faccf.zh-cn.xlf"

$deXlf1 = "39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.b0c490ba30254cd1a1156d1031ac43e4c04209be.de-de.xlf"
$deXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5ad68058f0fa4adb963ec49c3504f544c2eabaf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39ccbd6b-1989-4247-9b9b-ab8fe1dc01ed.b0c490ba30254cd1a1156d1031ac43e4c04209be.de-de.xlf"

$deXlf2 = "54d53a45-672b-45eb-a5f5-567c408ad3a1.10def667fa8a4b9195e9896236ea5151b05faccf.de-de.xlf"
$deXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5ad68058f0fa4adb963ec49c3504f544c2eabaf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/54d53a45-672b-45eb-a5f5-567c408ad3a1.10def667fa8a4b9195e9896236ea5151b05faccf.de-de.xlf"

# --- Overview sheet: refresh the Status column for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

$wsZh.Range("E2").Value = $md1
$wsZh.Range("F2").Value = $zhXlf1
$wsZh.Range("E3").Value = $md2
$wsZh.Range("F3").Value = $zhXlf2

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $md1Url, [Type]::Missing, [Type]::Missing, $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlf1Url, [Type]::Missing, [Type]::Missing, $zhXlf1)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $md2Url, [Type]::Missing, [Type]::Missing, $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlf2Url, [Type]::Missing, [Type]::Missing, $zhXlf2)

$wsZh.Range("E2:F3").Style = "Hyperlink"

$wsZh.Range("G2").Value = "2016-03-09 10:55:40"
$wsZh.Range("G3").Value = "2016-03-09 10:55:40"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

$wsDe.Range("E2").Value = $md1
$wsDe.Range("F2").Value = $deXlf1
$wsDe.Range("E3").Value = $md2
$wsDe.Range("F3").Value = $deXlf2

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $md1Url, [Type]::Missing, [Type]::Missing, $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlf1Url, [Type]::Missing, [Type]::Missing, $deXlf1)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $md2Url, [Type]::Missing, [Type]::Missing, $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlf2Url, [Type]::Missing, [Type]::Missing, $deXlf2)

$wsDe.Range("E2:F3").Style = "Hyperlink"

$wsDe.Range("G2").Value = "2016-03-09 10:55:51"
$wsDe.Range("G3").Value = "2016-03-09 10:55:51"
